$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "rotate animation" assertion text to the first empty row (A4).
# This introduces a new shared string and turns A4 into a text cell.
$ws.Range("A4").Value = "Mercedes-AMG Design: explore the possibilities."

# Column A is sized to best-fit its (longest) text; re-fit it now that A4 holds
# the new, longer string so the column widens like it does in a live Excel session
# (target best-fit width ~41.63 chars; feed a value that lands on the nearest
# width the host's column-width quantization can represent).
$ws.Columns("A").ColumnWidth = 40.833333333333336

# Move the active selection down to the next empty row (A5), as in the source session.
$ws.Range("A5").Select()
